$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 48.46865866666667
$ws.Range("H2").Value = 145.405976
$ws.Range("I2").Value = 0.1554430998624896
$ws.Range("J2").Value = 0.1554430998624896
$ws.Range("O2").Value = 0.6731047553411441
$ws.Range("P2").Value = 0.673104755341144
$ws.Range("Q2").Value = 10.09880047003023
$ws.Range("R2").Value = 90.88920423027201
$ws.Range("S2").Value = 0.1046294897024101
$ws.Range("T2").Value = 0.1046294897024101
$ws.Range("G3").Value = 48.46865866666667
$ws.Range("H3").Value = 145.405976
$ws.Range("I3").Value = 0.1554430998624896
$ws.Range("J3").Value = 0.1554430998624896
$ws.Range("M3").Value = 0.1011893333333333
$ws.Range("N3").Value = 0.303568
$ws.Range("O3").Value = 0.3268952446588559
$ws.Range("P3").Value = 0.3268952446588559
$ws.Range("Q3").Value = 4.904511258040889
$ws.Range("R3").Value = 44.14060132236801
$ws.Range("S3").Value = 0.05081361016007953
$ws.Range("T3").Value = 0.05081361016007953
$ws.Range("I4").Value = 0.20693808715897
$ws.Range("J4").Value = 0.20693808715897
$ws.Range("O4").Value = 0.6731047553411441
$ws.Range("P4").Value = 0.673104755341144
$ws.Range("S4").Value = 0.1392910105279029
$ws.Range("T4").Value = 0.1392910105279029
$ws.Range("I5").Value = 0.20693808715897
$ws.Range("J5").Value = 0.20693808715897
$ws.Range("M5").Value = 0.1011893333333333
$ws.Range("N5").Value = 0.303568
$ws.Range("O5").Value = 0.3268952446588559
$ws.Range("P5").Value = 0.3268952446588559
$ws.Range("Q5").Value = 6.529271348078222
$ws.Range("R5").Value = 58.763442132704
$ws.Range("S5").Value = 0.06764707663106716
$ws.Range("T5").Value = 0.06764707663106716
$ws.Range("G6").Value = 75.47903666666667
$ws.Range("H6").Value = 226.43711
$ws.Range("I6").Value = 0.2420676733554854
$ws.Range("J6").Value = 0.2420676733554854
$ws.Range("O6").Value = 0.6731047553411441
$ws.Range("P6").Value = 0.673104755341144
$ws.Range("Q6").Value = 15.72661080243556
$ws.Range("R6").Value = 141.53949722192
$ws.Range("S6").Value = 0.162936902049944
$ws.Range("T6").Value = 0.162936902049944
$ws.Range("G7").Value = 75.47903666666667
$ws.Range("H7").Value = 226.43711
$ws.Range("I7").Value = 0.2420676733554854
$ws.Range("J7").Value = 0.2420676733554854
$ws.Range("M7").Value = 0.1011893333333333
$ws.Range("N7").Value = 0.303568
$ws.Range("O7").Value = 0.3268952446588559
$ws.Range("P7").Value = 0.3268952446588559
$ws.Range("Q7").Value = 7.637673400942223
$ws.Range("R7").Value = 68.73906060848
$ws.Range("S7").Value = 0.07913077130554143
$ws.Range("T7").Value = 0.07913077130554143
$ws.Range("G8").Value = 51.18999233333333
$ws.Range("H8").Value = 153.569977
$ws.Range("I8").Value = 0.164170647777855
$ws.Range("J8").Value = 0.164170647777855
$ws.Range("O8").Value = 0.6731047553411441
$ws.Range("P8").Value = 0.673104755341144
$ws.Range("Q8").Value = 10.66581029592711
$ws.Range("R8").Value = 95.992292663344
$ws.Range("S8").Value = 0.1105040437067102
$ws.Range("T8").Value = 0.1105040437067102
$ws.Range("G9").Value = 51.18999233333333
$ws.Range("H9").Value = 153.569977
$ws.Range("I9").Value = 0.164170647777855
$ws.Range("J9").Value = 0.164170647777855
$ws.Range("M9").Value = 0.1011893333333333
$ws.Range("N9").Value = 0.303568
$ws.Range("O9").Value = 0.3268952446588559
$ws.Range("P9").Value = 0.3268952446588559
$ws.Range("Q9").Value = 5.179881197548444
$ws.Range("R9").Value = 46.618930777936
$ws.Range("S9").Value = 0.05366660407114477
$ws.Range("T9").Value = 0.05366660407114476
$ws.Range("G10").Value = 72.14667033333333
$ws.Range("H10").Value = 216.440011
$ws.Range("I10").Value = 0.2313804918452
$ws.Range("J10").Value = 0.2313804918452
$ws.Range("O10").Value = 0.6731047553411441
$ws.Range("P10").Value = 0.673104755341144
$ws.Range("Q10").Value = 15.03228783953245
$ws.Range("R10").Value = 135.290590555792
$ws.Range("S10").Value = 0.155743309354177
$ws.Range("T10").Value = 0.1557433093541769
$ws.Range("G11").Value = 72.14667033333333
$ws.Range("H11").Value = 216.440011
$ws.Range("I11").Value = 0.2313804918452
$ws.Range("J11").Value = 0.2313804918452
$ws.Range("M11").Value = 0.1011893333333333
$ws.Range("N11").Value = 0.303568
$ws.Range("O11").Value = 0.3268952446588559
$ws.Range("P11").Value = 0.3268952446588559
$ws.Range("Q11").Value = 7.300473473249778
$ws.Range("R11").Value = 65.704261259248
$ws.Range("S11").Value = 0.07563718249102309
$ws.Range("T11").Value = 0.07563718249102308
